$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(28, 9).Value = "[4]"
$ws.Cells.Item(49, 9).Value = "[4]"

$ws.Columns.Item(6).ColumnWidth = 53.85

$ws.Range("I49").Select()
